$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "244.88"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "24.95"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.101"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05642"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.500"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "2.977"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8102"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8355"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1326"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06972"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.02837"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09383"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001513"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0006023"
$ws.Range("E15").Value = "14OneONEWorstin24h"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006176"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.505"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.092"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3202"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.03171"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1324"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.744"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04684"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.1364"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001238"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004231"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.00009707"
$ws.Range("E27").Value = "26NitroExNTXBestin24h"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0001974"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03617"
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006285"
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1047"
$ws.Range("E42").Value = "41BKEXTokenBKK"
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002723"
$ws.Range("E43").Value = "42CEJICEJI"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007400"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005288"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000753"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.2008"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002286"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002108"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0002008"
